# Apply the "義大利2026.4" workbook update:
#  - Replace the old "tuscany or rome+firenze" / blank lodging rows (13-16)
#    on 工作表1 with "Tuscany" itinerary + "Lari" lodging (row 16 has no lodging).
#  - Update the saved selections on both sheets.
#  - Make 工作表1 the active/selected tab instead of 景點八選一.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 工作表1
$ws2 = $wb.Worksheets.Item(2)   # 景點八選一

# --- Update the itinerary / lodging cells on 工作表1 ------------------------
# ("Lari" is entered before "Tuscany" so the new shared-string entries land
#  in the same order as the source file: ... Torino, Lari, Tuscany)
$ws1.Range("C13").Value = "Lari"
$ws1.Range("B13").Value = "Tuscany"

$ws1.Range("C14").Value = "Lari"
$ws1.Range("B14").Value = "Tuscany"

$ws1.Range("C15").Value = "Lari"
$ws1.Range("B15").Value = "Tuscany"

$ws1.Range("B16").Value = "Tuscany"

# --- Update saved selections -------------------------------------------------
# 景點八選一 selection moves to G16 (done first so 工作表1 ends up active)
$ws2.Range("G16").Select()

# 工作表1 selection moves to B9, and 工作表1 becomes the active/selected tab
$ws1.Range("B9").Select()
